$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("flights") values are stored as text in most rows (inline string)
# but as real numbers in rows 11, 13-18. Preserve each cell's original type
# while updating the value, and strip any stray formatting (quote-prefix
# style) introduced by forcing a text entry so cell styles stay untouched.

$textUpdates = @{
    4  = "216"
    5  = "504"
    6  = "0"
    7  = "504"
    8  = "216"
    9  = "288"
    10 = "144"
    20 = "144"
    21 = "216"
    22 = "216"
    23 = "504"
    24 = "144"
    25 = "0"
    26 = "0"
    28 = "0"
    29 = "288"
    30 = "72"
    31 = "216"
    32 = "72"
    33 = "216"
    34 = "144"
    36 = "144"
    37 = "144"
    38 = "216"
    39 = "72"
    40 = "144"
    41 = "72"
    42 = "144"
    44 = "216"
    45 = "72"
    46 = "144"
    47 = "72"
    48 = "72"
    49 = "216"
    50 = "216"
    52 = "144"
    53 = "216"
    54 = "0"
    55 = "144"
    56 = "216"
    57 = "144"
    58 = "0"
    60 = "0"
    61 = "0"
    62 = "144"
    63 = "216"
    64 = "0"
    65 = "504"
    66 = "216"
}

$numberUpdates = @{
    11 = 0
    13 = 0
    14 = 72
    15 = 72
    16 = 144
    17 = 72
    18 = 72
}

foreach ($row in $textUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = "'" + $textUpdates[$row]
    $cell.ClearFormats()
}

foreach ($row in $numberUpdates.Keys) {
    $ws.Cells.Item($row, 6).Value = $numberUpdates[$row]
}
